$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.823.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.739.98'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '226.66'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.25%  '
$ws.Range("E6").Value = '  +0.04%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5158'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("E8").Value = '  +2.73%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '38.69'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -6.15%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06090'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.742.56'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.62%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07004'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.09%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '15.18'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6312'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.497'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '76.22'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.36%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '25.841.47'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.44'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.000006607'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.96%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.957.88'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.055'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.421'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.97%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.090'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.11%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '136.72'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.21%  '
$ws.Range("E27").Value = '  +3.81%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.817'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.25%  '
$ws.Range("E29").Value = '  -0.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '102.58'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08306'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.620'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.375'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.04404'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.608'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.58%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9677'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.10%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.5966'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.675'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01556'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.935'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.87%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9993'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.40'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.3803'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.39%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7240'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.876'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.05485'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '6.188'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1099'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '29.79'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.06%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.89'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.002'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.41%  '
